$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like numbers need to be forced to Text format
# so Excel stores them as the literal strings from the source data (matching
# the original inline-string / text cell contents) instead of converting them
# to numeric values.
$textCells = @('D4','D5','D7','D8','D9','D11','D12','D13','D14','D18','D19','D20','D22','D23','D24','D25','D26','D27','D28','D29','D31','D33','D35','D36','D37','D38','D39','D40','D41','D42','D43','D44','D45','D46','D48','D49','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.536.20'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '1.729.20'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '246.23'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D7').Value = '0.4831'
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').Value = '0.2665'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').Value = '0.06238'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').Value = '1.729.84'
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('D11').Value = '0.07087'
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('D12').Value = '15.59'
$ws.Range('D13').Value = '0.6085'
$ws.Range('E13').Value = '  -2.29%  '
$ws.Range('D14').Value = '4.552'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '26.535.25'
$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.000007254'
$ws.Range('E19').Value = '  +5.11%  '
$ws.Range('D20').Value = '11.51'
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('D21').Value = '1.955.98'
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('D22').Value = '4.506'
$ws.Range('E22').Value = '  -2.87%  '
$ws.Range('D23').Value = '8.767'
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('D24').Value = '5.241'
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('D25').Value = '137.25'
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').Value = '15.40'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').Value = '1.775'
$ws.Range('E27').Value = '  -2.43%  '
$ws.Range('D28').Value = '1.405'
$ws.Range('E28').Value = '  -2.14%  '
$ws.Range('D29').Value = '108.24'
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('E30').Value = '  -1.53%  '
$ws.Range('D31').Value = '0.08009'
$ws.Range('E31').Value = '  +1.69%  '
$ws.Range('E32').Value = '  -1.74%  '
$ws.Range('D33').Value = '0.04570'
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '2.619'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').Value = '0.6311'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('D38').Value = '0.8915'
$ws.Range('E38').Value = '  -5.73%  '
$ws.Range('D39').Value = '1.998'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').Value = '2.404'
$ws.Range('E40').Value = '  -0.77%  '
$ws.Range('D41').Value = '1.002'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('D42').Value = '0.01500'
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('D43').Value = '101.65'
$ws.Range('E43').Value = '  -10.18%  '
$ws.Range('D44').Value = '5.396'
$ws.Range('E44').Value = '  -6.26%  '
$ws.Range('D45').Value = '0.3891'
$ws.Range('E45').Value = '  -0.77%  '
$ws.Range('D46').Value = '7.040'
$ws.Range('E46').Value = '  +5.14%  '
$ws.Range('E47').Value = '  -2.23%  '
$ws.Range('D48').Value = '0.05395'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = '7.885'
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('D51').Value = '1.252'
$ws.Range('E51').Value = '  -1.59%  '
